$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value2 = "28.468.58"
$cell.Style = "Normal"
$ws.Range("E2").Value2 = "  +0.88%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value2 = "1.783.32"
$cell.Style = "Normal"
$ws.Range("E3").Value2 = "  -0.50%  "

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value2 = "1.007"
$cell.Style = "Normal"
$ws.Range("E4").Value2 = "  +0.28%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value2 = "328.07"
$cell.Style = "Normal"
$ws.Range("E5").Value2 = "  -2.78%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value2 = "1.002"
$cell.Style = "Normal"
$ws.Range("E6").Value2 = "  +0.20%  "

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value2 = "0.4386"
$cell.Style = "Normal"
$ws.Range("E7").Value2 = "  -3.00%  "

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value2 = "0.3742"
$cell.Style = "Normal"
$ws.Range("E8").Value2 = "  +5.19%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value2 = "45.11"
$cell.Style = "Normal"
$ws.Range("E9").Value2 = "  -0.80%  "

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value2 = "0.07598"
$cell.Style = "Normal"
$ws.Range("E10").Value2 = "  +1.55%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value2 = "1.134"
$cell.Style = "Normal"
$ws.Range("E11").Value2 = "  -0.39%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value2 = "22.50"
$cell.Style = "Normal"
$ws.Range("E12").Value2 = "  +0.51%  "

# Row 13
$ws.Range("E13").Value2 = "  +0.29%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value2 = "6.233"
$cell.Style = "Normal"
$ws.Range("E14").Value2 = "  +0.44%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value2 = "7.391"
$cell.Style = "Normal"
$ws.Range("E15").Value2 = "  +2.25%  "

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value2 = "1.789.62"
$cell.Style = "Normal"
$ws.Range("E16").Value2 = "  -0.04%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value2 = "0.00001089"
$cell.Style = "Normal"
$ws.Range("E17").Value2 = "  +0.45%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value2 = "0.06691"
$cell.Style = "Normal"
$ws.Range("E18").Value2 = "  +0.14%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value2 = "81.70"
$cell.Style = "Normal"
$ws.Range("E19").Value2 = "  +0.65%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value2 = "1.003"
$cell.Style = "Normal"
$ws.Range("E20").Value2 = "  +0.20%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value2 = "17.43"
$cell.Style = "Normal"
$ws.Range("E21").Value2 = "  +1.60%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value2 = "6.224"
$cell.Style = "Normal"
$ws.Range("E22").Value2 = "  -2.48%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value2 = "28.491.37"
$cell.Style = "Normal"
$ws.Range("E23").Value2 = "  +0.97%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value2 = "11.70"
$cell.Style = "Normal"
$ws.Range("E24").Value2 = "  -1.21%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value2 = "2.427"
$cell.Style = "Normal"
$ws.Range("E25").Value2 = "  +2.05%  "

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value2 = "20.46"
$cell.Style = "Normal"
$ws.Range("E26").Value2 = "  +0.11%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value2 = "152.56"
$cell.Style = "Normal"
$ws.Range("E27").Value2 = "  -0.70%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value2 = "2.358"
$cell.Style = "Normal"
$ws.Range("E28").Value2 = "  -0.61%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value2 = "1.994.84"
$cell.Style = "Normal"
$ws.Range("E29").Value2 = "  -0.04%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value2 = "1.301"
$cell.Style = "Normal"
$ws.Range("E30").Value2 = "  +2.50%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value2 = "129.29"
$cell.Style = "Normal"
$ws.Range("E31").Value2 = "  -2.26%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value2 = "3.969"
$cell.Style = "Normal"

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value2 = "5.852"
$cell.Style = "Normal"
$ws.Range("E33").Value2 = "  -0.18%  "

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value2 = "0.09216"
$cell.Style = "Normal"
$ws.Range("E34").Value2 = "  -2.05%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value2 = "0.2233"
$cell.Style = "Normal"
$ws.Range("E35").Value2 = "  +3.65%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value2 = "12.11"
$cell.Style = "Normal"
$ws.Range("E36").Value2 = "  +0.32%  "

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value2 = "0.6678"
$cell.Style = "Normal"
$ws.Range("E37").Value2 = "  +0.42%  "

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value2 = "0.06303"
$cell.Style = "Normal"
$ws.Range("E38").Value2 = "  +1.36%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value2 = "0.02333"
$cell.Style = "Normal"
$ws.Range("E39").Value2 = "  -1.42%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value2 = "5.210"
$cell.Style = "Normal"
$ws.Range("E40").Value2 = "  +0.85%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value2 = "1.201"
$cell.Style = "Normal"
$ws.Range("E41").Value2 = "  -0.56%  "

# Row 42
$ws.Range("E42").Value2 = "  -2.37%  "

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value2 = "8.047"
$cell.Style = "Normal"
$ws.Range("E43").Value2 = "  +0.02%  "

# Row 44
$ws.Range("B44").Value2 = "Frax"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value2 = "1.001"
$cell.Style = "Normal"
$ws.Range("E44").Value2 = "  +0.16%  "

# Row 45
$ws.Range("B45").Value2 = "EnergySwap"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value2 = "13.99"
$cell.Style = "Normal"
$ws.Range("E45").Value2 = "  +0.39%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value2 = "0.6104"
$cell.Style = "Normal"
$ws.Range("E46").Value2 = "  +0.81%  "

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value2 = "3.809"
$cell.Style = "Normal"
$ws.Range("E47").Value2 = "  -1.21%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value2 = "127.71"
$cell.Style = "Normal"
$ws.Range("E48").Value2 = "  -0.36%  "

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value2 = "2.017"
$cell.Style = "Normal"
$ws.Range("E49").Value2 = "  -0.13%  "

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value2 = "0.06979"
$cell.Style = "Normal"
$ws.Range("E50").Value2 = "  -1.49%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value2 = "1.132"
$cell.Style = "Normal"
$ws.Range("E51").Value2 = "  -2.47%  "
